$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.208.09'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').Value = '3.106.99'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''544.89'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').Value = '''141.62'
$ws.Range('E6').Value = '  +5.97%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.103.06'
$ws.Range('E8').Value = '  +2.12%  '
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('D10').Value = '''6.65'
$ws.Range('E10').Value = '  +4.52%  '
$ws.Range('D11').Value = '''0.158'
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  +7.08%  '
$ws.Range('D14').Value = '''35.16'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D15').Value = '3.616.23'
$ws.Range('E15').Value = '  +2.20%  '
$ws.Range('D16').Value = '64.249.58'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').Value = '3.111.09'
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '''485.71'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = '''13.50'
$ws.Range('E21').Value = '  +1.47%  '
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('D23').Value = '''7.17'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('D24').Value = '''79.58'
$ws.Range('E24').Value = '  +3.33%  '
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  +2.21%  '
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '''26.49'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('E32').Value = '  +3.80%  '
$ws.Range('D33').Value = '''57.99'
$ws.Range('E33').Value = '  -1.92%  '
$ws.Range('E34').Value = '  -3.46%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '''5.44'
$ws.Range('E35').Value = '  +7.75%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = '''499.49'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('D37').Value = '''6.07'
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('D38').Value = '3.289.77'
$ws.Range('E38').Value = '  +7.15%  '
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('D40').Value = '''0.0805'
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('E41').Value = '  +3.98%  '
$ws.Range('D42').Value = '''2.78'
$ws.Range('E42').Value = '  +7.97%  '
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('E44').Value = '  +2.47%  '
$ws.Range('D46').Value = '''124.31'
$ws.Range('E46').Value = '  +4.10%  '
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('D48').Value = '''25.38'
$ws.Range('E48').Value = '  +4.27%  '
$ws.Range('D49').Value = '0.0₃0539'
$ws.Range('E49').Value = '  +10.07%  '
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('D51').Value = '''2.42'
$ws.Range('E51').Value = '  +2.78%  '
